$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '328.65'
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = '1.13%'
$c.Style = "Normal"

$ws.Cells.Item(3, 2).Value = 'OKB'
$ws.Cells.Item(3, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '43.98'
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = '-1.34%'
$c.Style = "Normal"

$ws.Cells.Item(4, 2).Value = 'HuobiToken'
$ws.Cells.Item(4, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '5.517'
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = '0.58%'
$c.Style = "Normal"

$ws.Cells.Item(5, 2).Value = 'Cronos'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '0.08013'
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = '-0.56%'
$c.Style = "Normal"

$ws.Cells.Item(6, 2).Value = 'FTXToken'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '1.982'
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = '3.54%'
$c.Style = "Normal"

$ws.Cells.Item(7, 2).Value = 'GateToken'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '4.372'
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = '2.17%'
$c.Style = "Normal"

$ws.Cells.Item(8, 2).Value = 'BTSEToken'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '2.584'
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = '-4.82%'
$c.Style = "Normal"

$ws.Cells.Item(9, 2).Value = 'MXToken'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.9499'
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = '0.89%'
$c.Style = "Normal"

$ws.Cells.Item(10, 2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.1136'
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = '-2.41%'
$c.Style = "Normal"

$ws.Cells.Item(11, 2).Value = 'WazirX'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.1885'
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = '0.98%'
$c.Style = "Normal"

$ws.Cells.Item(12, 2).Value = 'MCDex'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '10.76'
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = '27.20%'
$c.Style = "Normal"

$ws.Cells.Item(13, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '0.09970'
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = '-1.06%'
$c.Style = "Normal"

$ws.Cells.Item(14, 2).Value = 'BitrueCoin'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '0.04778'
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = '13.31%'
$c.Style = "Normal"

$ws.Cells.Item(15, 2).Value = 'BitMartToken'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '0.1065'
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = '-0.16%'
$c.Style = "Normal"

$ws.Cells.Item(16, 2).Value = 'BitForexToken'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '0.001280'
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = '0.25%'
$c.Style = "Normal"

$ws.Cells.Item(17, 2).Value = 'CoinExToken'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '0.04072'
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = '-4.05%'
$c.Style = "Normal"

$ws.Cells.Item(18, 2).Value = 'TigerCash'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '0.005967'
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = '1.79%'
$c.Style = "Normal"

$ws.Cells.Item(19, 2).Value = 'LEO'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '3.364'
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = '-6.22%'
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = '-0.72%'
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '0.1416'
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = '2.77%'
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '0.2547'
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = '0.73%'
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '0.001265'
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = '2.29%'
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = '-4.91%'
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = '1.68%'
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '0.0003744'
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = '-6.19%'
$c.Style = "Normal"

$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = '-1.71%'
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '0.05653'
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = '3.18%'
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '0.007544'
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = '-1.53%'
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.1398'
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = '0.29%'
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.007411'
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = '3.20%'
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '0.002015'
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = '-0.43%'
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.008623'
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = '-0.46%'
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.00007108'
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = '-0.08%'
$c.Style = "Normal"

$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = '-0.06%'
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '0.003530'
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = '55.42%'
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '0.003777'
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = '6.96%'
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '0.00002100'
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = '-0.06%'
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.0002000'
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = '-0.06%'
$c.Style = "Normal"
